$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44533
$ws.Cells.Item(2, 12).Value = 'Primera'
$ws.Cells.Item(2, 14).Value = 18000
$ws.Cells.Item(2, 15).Value = 19000
$ws.Cells.Item(2, 16).Value = 18500
$ws.Cells.Item(2, 19).Value = 2312

# Row 3
$ws.Cells.Item(3, 4).Value = 44533
$ws.Cells.Item(3, 12).Value = 'Segunda'
$ws.Cells.Item(3, 13).Value = 100
$ws.Cells.Item(3, 14).Value = 16000
$ws.Cells.Item(3, 15).Value = 16000
$ws.Cells.Item(3, 16).Value = 16000
$ws.Cells.Item(3, 19).Value = 2000

# Row 4
$ws.Cells.Item(4, 4).Value = 44895
$ws.Cells.Item(4, 12).Value = 'Primera'
$ws.Cells.Item(4, 13).Value = 200
$ws.Cells.Item(4, 14).Value = 22000
$ws.Cells.Item(4, 15).Value = 22500
$ws.Cells.Item(4, 16).Value = 22250
$ws.Cells.Item(4, 19).Value = 2781

# Row 5
$ws.Cells.Item(5, 4).Value = 44495
$ws.Cells.Item(5, 12).Value = 'Segunda'
$ws.Cells.Item(5, 13).Value = 270
$ws.Cells.Item(5, 14).Value = 19000
$ws.Cells.Item(5, 15).Value = 20000
$ws.Cells.Item(5, 16).Value = 19556
$ws.Cells.Item(5, 19).Value = 2444

# Row 6
$ws.Cells.Item(6, 4).Value = 44876
$ws.Cells.Item(6, 12).Value = 'Primera'
$ws.Cells.Item(6, 13).Value = 300
$ws.Cells.Item(6, 14).Value = 22000
$ws.Cells.Item(6, 15).Value = 22500
$ws.Cells.Item(6, 16).Value = 22250
$ws.Cells.Item(6, 19).Value = 2781

# Row 7
$ws.Cells.Item(7, 4).Value = 44519
$ws.Cells.Item(7, 12).Value = 'Primera'
$ws.Cells.Item(7, 13).Value = 400
$ws.Cells.Item(7, 14).Value = 21000
$ws.Cells.Item(7, 15).Value = 22000
$ws.Cells.Item(7, 16).Value = 21500
$ws.Cells.Item(7, 19).Value = 2688

# Row 8
$ws.Cells.Item(8, 4).Value = 44519
$ws.Cells.Item(8, 12).Value = 'Segunda'
$ws.Cells.Item(8, 14).Value = 18000
$ws.Cells.Item(8, 15).Value = 18000
$ws.Cells.Item(8, 16).Value = 18000
$ws.Cells.Item(8, 19).Value = 2250

# Row 9
$ws.Cells.Item(9, 4).Value = 44880

# Row 10
$ws.Cells.Item(10, 4).Value = 44523
$ws.Cells.Item(10, 13).Value = 400
$ws.Cells.Item(10, 15).Value = 22000
$ws.Cells.Item(10, 16).Value = 21500
$ws.Cells.Item(10, 19).Value = 2688

# Row 11
$ws.Cells.Item(11, 4).Value = 44523
$ws.Cells.Item(11, 13).Value = 100
$ws.Cells.Item(11, 15).Value = 18000
$ws.Cells.Item(11, 16).Value = 18000
$ws.Cells.Item(11, 19).Value = 2250

# Row 12
$ws.Cells.Item(12, 4).Value = 44498
$ws.Cells.Item(12, 12).Value = 'Segunda'
$ws.Cells.Item(12, 14).Value = 19000
$ws.Cells.Item(12, 15).Value = 20000
$ws.Cells.Item(12, 16).Value = 19500
$ws.Cells.Item(12, 19).Value = 2438

# Row 13
$ws.Cells.Item(13, 4).Value = 44873
$ws.Cells.Item(13, 13).Value = 300
$ws.Cells.Item(13, 14).Value = 22000
$ws.Cells.Item(13, 15).Value = 22500
$ws.Cells.Item(13, 16).Value = 22250
$ws.Cells.Item(13, 19).Value = 2781

# Row 14
$ws.Cells.Item(14, 4).Value = 44488
$ws.Cells.Item(14, 13).Value = 160
$ws.Cells.Item(14, 14).Value = 17000
$ws.Cells.Item(14, 16).Value = 17500
$ws.Cells.Item(14, 19).Value = 2188

# Row 16
$ws.Cells.Item(16, 4).Value = 44530
$ws.Cells.Item(16, 12).Value = 'Primera'
$ws.Cells.Item(16, 13).Value = 200

# Row 17
$ws.Cells.Item(17, 4).Value = 44530
$ws.Cells.Item(17, 12).Value = 'Segunda'
$ws.Cells.Item(17, 13).Value = 100
$ws.Cells.Item(17, 14).Value = 16000
$ws.Cells.Item(17, 15).Value = 16000
$ws.Cells.Item(17, 16).Value = 16000
$ws.Cells.Item(17, 17).Value = '$/bandeja 8 kilos'
$ws.Cells.Item(17, 19).Value = 2000
$ws.Cells.Item(17, 20).Value = 8

# Row 18
$ws.Cells.Item(18, 4).Value = 44894
$ws.Cells.Item(18, 13).Value = 200
$ws.Cells.Item(18, 14).Value = 22000
$ws.Cells.Item(18, 15).Value = 22500
$ws.Cells.Item(18, 16).Value = 22250
$ws.Cells.Item(18, 19).Value = 2781

# Row 19
$ws.Cells.Item(19, 4).Value = 44890
$ws.Cells.Item(19, 12).Value = 'Primera'
$ws.Cells.Item(19, 13).Value = 200
$ws.Cells.Item(19, 14).Value = 22000
$ws.Cells.Item(19, 15).Value = 22500
$ws.Cells.Item(19, 16).Value = 22250
$ws.Cells.Item(19, 19).Value = 2781

# Row 20
$ws.Cells.Item(20, 4).Value = 44509
$ws.Cells.Item(20, 13).Value = 200
$ws.Cells.Item(20, 14).Value = 19000
$ws.Cells.Item(20, 15).Value = 20000
$ws.Cells.Item(20, 16).Value = 19500
$ws.Cells.Item(20, 19).Value = 2438

# Row 21
$ws.Cells.Item(21, 4).Value = 44516
$ws.Cells.Item(21, 14).Value = 18000
$ws.Cells.Item(21, 15).Value = 19000
$ws.Cells.Item(21, 16).Value = 18500
$ws.Cells.Item(21, 19).Value = 2312

# Row 22
$ws.Cells.Item(22, 4).Value = 44512
$ws.Cells.Item(22, 12).Value = 'Segunda'
$ws.Cells.Item(22, 14).Value = 19000
$ws.Cells.Item(22, 15).Value = 20000
$ws.Cells.Item(22, 16).Value = 19500
$ws.Cells.Item(22, 19).Value = 2438

# Row 23
$ws.Cells.Item(23, 4).Value = 44159
$ws.Cells.Item(23, 13).Value = 300
$ws.Cells.Item(23, 14).Value = 2000
$ws.Cells.Item(23, 15).Value = 2100
$ws.Cells.Item(23, 16).Value = 2050
$ws.Cells.Item(23, 17).Value = '$/kilo (en caja de 14 kilos)'
$ws.Cells.Item(23, 19).Value = 2050
$ws.Cells.Item(23, 20).Value = 1

# Row 24
$ws.Cells.Item(24, 4).Value = 44491
$ws.Cells.Item(24, 12).Value = 'Segunda'
$ws.Cells.Item(24, 14).Value = 18000
$ws.Cells.Item(24, 15).Value = 19000
$ws.Cells.Item(24, 16).Value = 18500
$ws.Cells.Item(24, 19).Value = 2312

# Row 25
$ws.Cells.Item(25, 4).Value = 44526
$ws.Cells.Item(25, 12).Value = 'Primera'
$ws.Cells.Item(25, 13).Value = 300
$ws.Cells.Item(25, 14).Value = 21000
$ws.Cells.Item(25, 15).Value = 21000
$ws.Cells.Item(25, 16).Value = 21000
$ws.Cells.Item(25, 19).Value = 2625

# Row 26
$ws.Cells.Item(26, 4).Value = 44505
$ws.Cells.Item(26, 13).Value = 300
$ws.Cells.Item(26, 16).Value = 19500
$ws.Cells.Item(26, 19).Value = 2438
